$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '80.574.29'
$ws.Range("E2").Value = '  +5.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.174.07'
$ws.Range("E3").Value = '  +3.05%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.87'
$ws.Range("E5").Value = '  +5.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '623.57'
$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.272'
$ws.Range("E7").Value = '  +26.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.587'
$ws.Range("E9").Value = '  +6.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.177.94'
$ws.Range("E10").Value = '  +3.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.584'
$ws.Range("E11").Value = '  +27.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000255'
$ws.Range("E12").Value = '  +25.83%  '

$ws.Range("E13").Value = '  +1.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.788.51'
$ws.Range("E14").Value = '  +3.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.24'
$ws.Range("E15").Value = '  -0.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.69'
$ws.Range("E16").Value = '  +7.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.715.86'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.187.08'
$ws.Range("E18").Value = '  +3.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.20'
$ws.Range("E19").Value = '  +4.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.00'
$ws.Range("E20").Value = '  +9.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.16'
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '436.50'
$ws.Range("E22").Value = '  +12.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.15'
$ws.Range("E23").Value = '  +13.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.96'
$ws.Range("E24").Value = '  +8.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '76.03'
$ws.Range("E26").Value = '  +4.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.67'
$ws.Range("E27").Value = '  +1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.84'
$ws.Range("E28").Value = '  +4.08%  '

$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000121'
$ws.Range("E30").Value = '  +8.55%  '

$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.92'
$ws.Range("E32").Value = '  +6.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '556.86'
$ws.Range("E33").Value = '  +9.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").Value = '  +1.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.151'
$ws.Range("E35").Value = '  +14.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.98'
$ws.Range("E36").Value = '  +2.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.93'
$ws.Range("E37").Value = '  +9.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.123'
$ws.Range("E38").Value = '  +19.85%  '

$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.403'
$ws.Range("E40").Value = '  +6.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '20.80'
$ws.Range("E41").Value = '  +3.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '164.48'
$ws.Range("E42").Value = '  +0.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.64'
$ws.Range("E43").Value = '  +7.37%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '189.90'
$ws.Range("E45").Value = '  -2.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.80'
$ws.Range("E46").Value = '  +7.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.68'
$ws.Range("E47").Value = '  +7.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.779'
$ws.Range("E48").Value = '  -2.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("E49").Value = '  +1.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.84'
$ws.Range("E50").Value = '  +4.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.24'
$ws.Range("E51").Value = '  +7.29%  '
